$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 379.30768
$ws.Range("I80").Value = 392.875
$ws.Range("J80").Value = 357.6
$ws.Range("K80").Value = 1178.625
$ws.Range("L80").Value = 1072.8
$ws.Range("M80").Value = -180.625
$ws.Range("N80").Value = -3068.8
$ws.Range("H83").Value = 379.30768
$ws.Range("I83").Value = 392.875
$ws.Range("J83").Value = 357.6
$ws.Range("K83").Value = 3535.875
$ws.Range("L83").Value = 3218.4
$ws.Range("M83").Value = 1456.125
$ws.Range("N83").Value = -13202.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1059
$ws.Range("I32").Value = 1059
$ws.Range("K32").Value = 1059
$ws.Range("M32").Value = -772
$ws.Range("H45").Value = 2435.2856
$ws.Range("I45").Value = 1391.4286
$ws.Range("K45").Value = 1391.4286
$ws.Range("M45").Value = -1014.4286
$ws.Range("H61").Value = 6248.625
$ws.Range("I61").Value = 5331.5
$ws.Range("K61").Value = 5331.5
$ws.Range("M61").Value = -5119.5
$ws.Range("H80").Value = 36264
$ws.Range("I80").Value = 16550
$ws.Range("J80").Value = 49406.668
$ws.Range("K80").Value = 16550
$ws.Range("L80").Value = 49406.668
$ws.Range("M80").Value = -15552
$ws.Range("N80").Value = -51402.668
$ws.Range("H83").Value = 36264
$ws.Range("I83").Value = 16550
$ws.Range("J83").Value = 49406.668
$ws.Range("K83").Value = 49650
$ws.Range("L83").Value = 148220.004
$ws.Range("M83").Value = -44658
$ws.Range("N83").Value = -158204.004
$ws.Range("H132").Value = 2496
$ws.Range("I132").Value = 2496
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7488
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -4958
$ws.Range("H136").Value = 6248.625
$ws.Range("I136").Value = 5331.5
$ws.Range("K136").Value = 15994.5
$ws.Range("M136").Value = -13444.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2709.2778
$ws.Range("I20").Value = 1622.6666
$ws.Range("J20").Value = 3252.5833
$ws.Range("K20").Value = 1622.6666
$ws.Range("L20").Value = 3252.5833
$ws.Range("M20").Value = -1375.6666
$ws.Range("N20").Value = -3746.5833
$ws.Range("H82").Value = 38833.625
$ws.Range("I82").Value = 15124.5
$ws.Range("K82").Value = 15124.5
$ws.Range("M82").Value = -14741.5
$ws.Range("H85").Value = 38833.625
$ws.Range("I85").Value = 15124.5
$ws.Range("K85").Value = 15124.5
$ws.Range("M85").Value = -13798.5
$ws.Range("H105").Value = 1000.86365
$ws.Range("I105").Value = 958.7222
$ws.Range("J105").Value = 1190.5
$ws.Range("K105").Value = 958.7222
$ws.Range("L105").Value = 1190.5
$ws.Range("M105").Value = 788.2778
$ws.Range("N105").Value = -4684.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5555.5713
$ws.Range("I58").Value = 4978.8
$ws.Range("K58").Value = 4978.8
$ws.Range("M58").Value = -4775.8
$ws.Range("H103").Value = 11909.667
$ws.Range("I103").Value = 11909.667
$ws.Range("K103").Value = 11909.667
$ws.Range("M103").Value = -10737.667
$ws.Range("H132").Value = 2050.375
$ws.Range("I132").Value = 2050.375
$ws.Range("K132").Value = 6151.125
$ws.Range("M132").Value = -3621.125
$ws.Range("H136").Value = 5555.5713
$ws.Range("I136").Value = 4978.8
$ws.Range("K136").Value = 14936.4
$ws.Range("M136").Value = -12386.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 960.5625
$ws.Range("J5").Value = 579.6
$ws.Range("L5").Value = 1738.8
$ws.Range("N5").Value = -1962.8
$ws.Range("H6").Value = 10
$ws.Range("I6").Value = 10
$ws.Range("K6").Value = 30
$ws.Range("M6").Value = 83
$ws.Range("H64").Value = 1600
$ws.Range("I64").Value = 1600
$ws.Range("K64").Value = 4800
$ws.Range("M64").Value = -4530
$ws.Range("H67").Value = 1600
$ws.Range("I67").Value = 1600
$ws.Range("K67").Value = 4800
$ws.Range("M67").Value = -3864
$ws.Range("H135").Value = 960.5625
$ws.Range("J135").Value = 579.6
$ws.Range("L135").Value = 5216.400000000001
$ws.Range("N135").Value = -10286.4
$ws.Range("H140").Value = 2056.4546
$ws.Range("I140").Value = 1762.1
$ws.Range("K140").Value = 5286.299999999999
$ws.Range("M140").Value = -106.2999999999993

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10008
$ws.Range("I70").Value = 10008
$ws.Range("K70").Value = 10008
$ws.Range("M70").Value = -9738
$ws.Range("H73").Value = 10008
$ws.Range("I73").Value = 10008
$ws.Range("K73").Value = 10008
$ws.Range("M73").Value = -9072
$ws.Range("H102").Value = 1483.7368
$ws.Range("I102").Value = 1113.4286
$ws.Range("J102").Value = 2520.6
$ws.Range("K102").Value = 1113.4286
$ws.Range("L102").Value = 2520.6
$ws.Range("M102").Value = 508.5714
$ws.Range("N102").Value = -5764.6
$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("M108").ClearContents()
$ws.Range("H126").Value = 1930.5
$ws.Range("I126").Value = 1833.3334
$ws.Range("J126").Value = 2222
$ws.Range("K126").Value = 5500.0002
$ws.Range("L126").Value = 6666
$ws.Range("M126").Value = -3030.0002
$ws.Range("N126").Value = -11606
$ws.Range("H132").Value = 3019.8572
$ws.Range("I132").Value = 3190.8333
$ws.Range("J132").Value = 1994
$ws.Range("K132").Value = 9572.499899999999
$ws.Range("L132").Value = 5982
$ws.Range("M132").Value = -7042.499899999999
$ws.Range("N132").Value = -11042

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8000
$ws.Range("I7").Value = 8000
$ws.Range("K7").Value = 8000
$ws.Range("M7").Value = -7888
$ws.Range("H61").Value = 5899.5557
$ws.Range("I61").Value = 3274.5
$ws.Range("K61").Value = 3274.5
$ws.Range("M61").Value = -3072.5
$ws.Range("H74").Value = 30000
$ws.Range("J74").Value = 50000
$ws.Range("L74").Value = 50000
$ws.Range("N74").Value = -51996
$ws.Range("H77").Value = 30000
$ws.Range("J77").Value = 50000
$ws.Range("L77").Value = 150000
$ws.Range("N77").Value = -159984
$ws.Range("H100").Value = 5840.5
$ws.Range("I100").Value = 1217.2
$ws.Range("J100").Value = 9142.857
$ws.Range("K100").Value = 1217.2
$ws.Range("L100").Value = 9142.857
$ws.Range("M100").Value = -676.2
$ws.Range("N100").Value = -10224.857
$ws.Range("H113").Value = 5899.5557
$ws.Range("I113").Value = 3274.5
$ws.Range("K113").Value = 3274.5
$ws.Range("M113").Value = -1104.5
$ws.Range("H126").Value = 8000
$ws.Range("I126").Value = 8000
$ws.Range("K126").Value = 24000
$ws.Range("M126").Value = -21530
$ws.Range("H132").Value = 5266.3
$ws.Range("I132").Value = 5096.2856
$ws.Range("K132").Value = 15288.8568
$ws.Range("M132").Value = -12758.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 921.7368
$ws.Range("I100").Value = 758.0714
$ws.Range("J100").Value = 1380
$ws.Range("K100").Value = 1516.1428
$ws.Range("L100").Value = 2760
$ws.Range("M100").Value = -975.1428000000001
$ws.Range("N100").Value = -3842
$ws.Range("H113").Value = 457.85715
$ws.Range("I113").Value = 319
$ws.Range("K113").Value = 957
$ws.Range("M113").Value = 1213
$ws.Range("H136").Value = 3527.087
$ws.Range("I136").Value = 2490.9285
$ws.Range("K136").Value = 7472.7855
$ws.Range("M136").Value = -4922.7855
